# Atualização de bases das ligas, do dia: 19-06-2024 às 21:51
# Swap the full data (all columns except A "index" and D "Date") between
# the given pairs of rows. This corrects mismatched row assignments in the
# source feed (ids/teams/odds swapped between two fixtures on the same date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$row1,
        [int]$row2,
        [int]$firstCol = 2,   # column B
        [int]$lastCol = 30    # column AD
    )

    $vals1 = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals1 += ,$ws.Cells.Item($row1, $c).Value2
    }
    $vals2 = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals2 += ,$ws.Cells.Item($row2, $c).Value2
    }

    for ($i = 0; $i -lt $vals1.Count; $i++) {
        $c = $firstCol + $i
        $ws.Cells.Item($row1, $c).Value2 = $vals2[$i]
        $ws.Cells.Item($row2, $c).Value2 = $vals1[$i]
    }
}

Swap-RowData -row1 130 -row2 131
Swap-RowData -row1 134 -row2 137
Swap-RowData -row1 135 -row2 136
Swap-RowData -row1 139 -row2 140
Swap-RowData -row1 254 -row2 256
